$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.293.82'
$ws.Range('E2').Value = '  -1.08%  '
$ws.Range('D3').Value = '2.253.05'
$ws.Range('E3').Value = '  -1.07%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').Value = "'247.66"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.49%  '
$ws.Range('D6').Value = "'0.621"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.52%  '
$ws.Range('D7').Value = "'74.44"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.16%  '
$ws.Range('D9').Value = "'0.613"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.63%  '
$ws.Range('D10').Value = "'41.46"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +4.52%  '
$ws.Range('E11').Value = '  -3.01%  '
$ws.Range('D12').Value = "'7.13"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -4.58%  '
$ws.Range('E13').Value = '  -2.85%  '
$ws.Range('D14').Value = '2.591.89'
$ws.Range('E14').Value = '  -1.04%  '
$ws.Range('D15').Value = "'14.56"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.32%  '
$ws.Range('E16').Value = '  -1.58%  '
$ws.Range('D17').Value = '2.262.14'
$ws.Range('E17').Value = '  -0.85%  '
$ws.Range('D18').Value = '42.191.65'
$ws.Range('E18').Value = '  -1.11%  '
$ws.Range('E19').Value = '  -1.73%  '
$ws.Range('D20').Value = "'6.12"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.59%  '
$ws.Range('D21').Value = "'71.87"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.94%  '
$ws.Range('E22').Value = '  +4.53%  '
$ws.Range('D23').Value = "'232.04"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.07%  '
$ws.Range('E24').Value = '  +0.05%  '
$ws.Range('D25').Value = "'11.19"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.24%  '
$ws.Range('D26').Value = "'7.90"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +24.41%  '
$ws.Range('D27').Value = "'3.55"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -7.90%  '
$ws.Range('E28').Value = '  -3.91%  '
$ws.Range('D29').Value = "'2.17"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.27%  '
$ws.Range('D30').Value = "'168.99"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.54%  '
$ws.Range('D31').Value = "'20.72"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.50%  '
$ws.Range('D32').Value = "'0.0824"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -5.92%  '
$ws.Range('B33').Value = 'InjectiveProtocol'
$ws.Range('C33').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D33').Value = "'30.93"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.25%  '
$ws.Range('B34').Value = 'Kaspa'
$ws.Range('C34').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D34').Value = "'0.119"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -6.05%  '
$ws.Range('E35').Value = '  -2.30%  '
$ws.Range('D36').Value = "'4.52"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.68%  '
$ws.Range('D37').Value = "'4.89"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.59%  '
$ws.Range('D38').Value = "'0.0307"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.23%  '
$ws.Range('D39').Value = "'13.53"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.52%  '
$ws.Range('D40').Value = "'2.18"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -5.15%  '
$ws.Range('D41').Value = "'5.79"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.32%  '
$ws.Range('D42').Value = "'61.67"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.36%  '
$ws.Range('D43').Value = "'0.203"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.32%  '
$ws.Range('D44').Value = "'108.08"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.47%  '
$ws.Range('D45').Value = "'8.67"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.50%  '
$ws.Range('E47').Value = '  -0.29%  '
$ws.Range('D48').Value = "'1.12"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.43%  '
$ws.Range('E49').Value = '  -0.48%  '
$ws.Range('D50').Value = "'2.27"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.08%  '
$ws.Range('D51').Value = "'4.13"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.01%  '
